$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.59145724773407
$ws.Range("B1").Value = 2.467047691345215
$ws.Range("C1").Value = 2.178581714630127
$ws.Range("D1").Value = 1.813200831413269
$ws.Range("E1").Value = 1.718773603439331
